# Add "NA" values under the duplicate_image_filename column (E) for rows 2-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
